$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update probability/value column (D) for several conditions (prior distributions)
$ws.Range("D2").Value2 = 0.05
$ws.Range("D3").Value2 = 0.05
$ws.Range("D6").Value2 = 0.22
$ws.Range("D7").Value2 = 0.22
$ws.Range("D8").Value2 = 0.47
$ws.Range("D9").Value2 = 0.47
$ws.Range("D10").Value2 = 1
$ws.Range("D11").Value2 = 1

# Remove the two extra "pop-out" levels (old rows 12:15 -- cont6_same/opp,
# cont7_same/opp), which shifts the final "blank" row (old row 16) up into
# row 12 -- this keeps row 12's style for column B intact.
$ws.Rows("12:15").Delete()

# Row 12 (after the shift) inherited column A from the old row 16 (15);
# restore it to the original row 12's condition number (11), which is the
# only column that must NOT come from the old row 16.
$ws.Range("A12").Value2 = 11

# Reset the active selection to the new last data cell.
$ws.Range("D12").Select() | Out-Null
